# PowerShell COM-interop script to apply the gh-pages data refresh diff
# to the '上海-漫展信息' workbook (4 sheets: 展览/演出/本地生活/全部类型).

$wb = $excel.ActiveWorkbook

# --- Sheet '展览': refresh 'want to go' counts (column F) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4,6).Value = 358
$ws.Cells.Item(5,6).Value = 1592
$ws.Cells.Item(6,6).Value = 774
$ws.Cells.Item(7,6).Value = 688
$ws.Cells.Item(8,6).Value = 1282
$ws.Cells.Item(9,6).Value = 2558
$ws.Cells.Item(10,6).Value = 1335
$ws.Cells.Item(11,6).Value = 302
$ws.Cells.Item(12,6).Value = 2305
$ws.Cells.Item(13,6).Value = 2009
$ws.Cells.Item(14,6).Value = 709
$ws.Cells.Item(15,6).Value = 6246
$ws.Cells.Item(16,6).Value = 115
$ws.Cells.Item(17,6).Value = 1212
$ws.Cells.Item(18,6).Value = 133
$ws.Cells.Item(19,6).Value = 1426
$ws.Cells.Item(20,6).Value = 1319
$ws.Cells.Item(22,6).Value = 98
$ws.Cells.Item(23,6).Value = 2128
$ws.Cells.Item(25,6).Value = 688
$ws.Cells.Item(26,6).Value = 219
$ws.Cells.Item(27,6).Value = 0
$ws.Cells.Item(28,6).Value = 277
$ws.Cells.Item(29,6).Value = 1241
$ws.Cells.Item(31,6).Value = 3664
$ws.Cells.Item(32,6).Value = 634
$ws.Cells.Item(33,6).Value = 1662
$ws.Cells.Item(34,6).Value = 502
$ws.Cells.Item(35,6).Value = 147
$ws.Cells.Item(38,6).Value = 455
$ws.Cells.Item(39,6).Value = 377
$ws.Cells.Item(40,6).Value = 1762
$ws.Cells.Item(41,6).Value = 39
$ws.Cells.Item(42,6).Value = 99
$ws.Cells.Item(43,6).Value = 885
$ws.Cells.Item(44,6).Value = 1043
$ws.Cells.Item(45,6).Value = 506
$ws.Cells.Item(49,6).Value = 70

# --- Sheet '演出': refresh 'want to go' counts (column F) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3,6).Value = 77
$ws.Cells.Item(7,6).Value = 437
$ws.Cells.Item(11,6).Value = 380
$ws.Cells.Item(21,6).Value = 581
$ws.Cells.Item(22,6).Value = 230
$ws.Cells.Item(23,6).Value = 351
$ws.Cells.Item(26,6).Value = 79
$ws.Cells.Item(27,6).Value = 79
$ws.Cells.Item(30,6).Value = 299
$ws.Cells.Item(32,6).Value = 133
$ws.Cells.Item(35,6).Value = 40
$ws.Cells.Item(38,6).Value = 183

# --- Sheet '本地生活': refresh 'want to go' counts (column F) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4,6).Value = 3290
$ws.Cells.Item(5,6).Value = 393
$ws.Cells.Item(7,6).Value = 1446
$ws.Cells.Item(9,6).Value = 379
$ws.Cells.Item(10,6).Value = 2744
$ws.Cells.Item(11,6).Value = 263
$ws.Cells.Item(12,6).Value = 495
$ws.Cells.Item(13,6).Value = 464
$ws.Cells.Item(14,6).Value = 1125

# --- Sheet '全部类型': refresh 'want to go' counts (column F) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2,6).Value = 1446
$ws.Cells.Item(4,6).Value = 358
$ws.Cells.Item(5,6).Value = 379
$ws.Cells.Item(6,6).Value = 2744
$ws.Cells.Item(7,6).Value = 1592
$ws.Cells.Item(8,6).Value = 774
$ws.Cells.Item(9,6).Value = 688
$ws.Cells.Item(10,6).Value = 1282
$ws.Cells.Item(11,6).Value = 2558
$ws.Cells.Item(12,6).Value = 1335
$ws.Cells.Item(22,6).Value = 1427
$ws.Cells.Item(23,6).Value = 1319
$ws.Cells.Item(25,6).Value = 2128
$ws.Cells.Item(26,6).Value = 351
$ws.Cells.Item(27,6).Value = 79
$ws.Cells.Item(29,6).Value = 688
$ws.Cells.Item(30,6).Value = 219
$ws.Cells.Item(31,6).Value = 5239
$ws.Cells.Item(32,6).Value = 277
$ws.Cells.Item(33,6).Value = 1241
$ws.Cells.Item(34,6).Value = 3664
$ws.Cells.Item(35,6).Value = 299
$ws.Cells.Item(36,6).Value = 1662
$ws.Cells.Item(37,6).Value = 505
$ws.Cells.Item(38,6).Value = 147
$ws.Cells.Item(40,6).Value = 377
$ws.Cells.Item(41,6).Value = 1762
$ws.Cells.Item(42,6).Value = 39
$ws.Cells.Item(43,6).Value = 40
$ws.Cells.Item(44,6).Value = 99
$ws.Cells.Item(45,6).Value = 885
$ws.Cells.Item(46,6).Value = 1043
$ws.Cells.Item(47,6).Value = 506
$ws.Cells.Item(48,6).Value = 183
$ws.Cells.Item(49,6).Value = 183
$ws.Cells.Item(51,6).Value = 70

# --- Sheet '演出': row 12 ticket price -> now marked unavailable ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(12,7).Value = '不可售'

# --- Sheet '全部类型': row 13 event removed (now unsellable) from the merged,
#     date-sorted listing -> rows 13-20 shift up one, new row 21 appended ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(13,2).Value = '2024-07-20'
$ws.Cells.Item(13,3).Value = '上海·漫游L+动漫游戏嘉年华（免费展）'
$ws.Cells.Item(13,4).Value = '申长路869号 上海龙湖虹桥天街'
$ws.Cells.Item(13,5).Value = '2024.07.20 10:00-07.21 17:00'
$ws.Cells.Item(13,6).Value = 302
$ws.Cells.Item(13,7).Value = 20
$ws.Cells.Item(13,8).Value = 'https://show.bilibili.com/platform/detail.html?id=88134'
$ws.Cells.Item(13,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/adaE6Z6f1719454819535.jpeg'

$ws.Cells.Item(14,2).Value = '2024-07-20'
$ws.Cells.Item(14,3).Value = '上海·第九届Redamancy动漫游戏嘉年华'
$ws.Cells.Item(14,4).Value = '中山北路3300号4楼 上海环球港'
$ws.Cells.Item(14,6).Value = 2305
$ws.Cells.Item(14,7).Value = 60
$ws.Cells.Item(14,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84637'
$ws.Cells.Item(14,9).Value = '//i1.hdslb.com/bfs/openplatform/202404/hWLkXqwM1713194236349.png'

$ws.Cells.Item(15,2).Value = '2024-07-21'
$ws.Cells.Item(15,3).Value = '上海·原神×星穹铁道ONLY 2.0'
$ws.Cells.Item(15,4).Value = '吴中路1588号上海爱琴海购物中心F4 竞梦元宇宙'
$ws.Cells.Item(15,5).Value = '2024.07.21 10:00-07.21 17:00'
$ws.Cells.Item(15,6).Value = 2009
$ws.Cells.Item(15,7).Value = 68
$ws.Cells.Item(15,8).Value = 'https://show.bilibili.com/platform/detail.html?id=88273'
$ws.Cells.Item(15,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/OPX4neRo1719567526505.png'

$ws.Cells.Item(16,2).Value = '2024-07-21'
$ws.Cells.Item(16,3).Value = '上海·葬送的芙莉莲ONLY'
$ws.Cells.Item(16,4).Value = '逸仙路301号靠纪念路路口 上海宝丰联大酒店'
$ws.Cells.Item(16,6).Value = 709
$ws.Cells.Item(16,7).Value = 65
$ws.Cells.Item(16,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85193'
$ws.Cells.Item(16,9).Value = '//i2.hdslb.com/bfs/openplatform/202404/VIM2lfxY1714361685906.jpeg'

$ws.Cells.Item(17,2).Value = '2024-07-26'
$ws.Cells.Item(17,3).Value = '上海·2024ChinaJoy中国国际数码互动娱乐展览会 '
$ws.Cells.Item(17,4).Value = '龙阳路2345号 上海新国际博览中心'
$ws.Cells.Item(17,5).Value = '2024.07.26 09:00-07.29 16:00'
$ws.Cells.Item(17,6).Value = 6246
$ws.Cells.Item(17,7).Value = 100
$ws.Cells.Item(17,8).Value = 'https://show.bilibili.com/platform/detail.html?id=88037'
$ws.Cells.Item(17,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/0yTYZsPC1719306558643.jpeg'

$ws.Cells.Item(18,2).Value = '2024-07-26'
$ws.Cells.Item(18,3).Value = '上海·BH夏日欢愉pro '
$ws.Cells.Item(18,4).Value = '人民大道221号 迪美购物中心'
$ws.Cells.Item(18,5).Value = '2024.07.26 10:30-07.28 18:00'
$ws.Cells.Item(18,6).Value = 115
$ws.Cells.Item(18,7).Value = 29.9
$ws.Cells.Item(18,8).Value = 'https://show.bilibili.com/platform/detail.html?id=88604'
$ws.Cells.Item(18,9).Value = '//i2.hdslb.com/bfs/openplatform/202407/5ArdQNaN1719995206315.png'

$ws.Cells.Item(19,2).Value = '2024-07-26'
$ws.Cells.Item(19,3).Value = '上海·盗墓笔记官方授权「四季同书」主题店'
$ws.Cells.Item(19,4).Value = '南京东路830号第一百货商业中心B馆5楼(海底捞旁边) 第一百货商业中心'
$ws.Cells.Item(19,5).Value = '2024.07.26 00:00-09.28 23:59'
$ws.Cells.Item(19,6).Value = 495
$ws.Cells.Item(19,7).Value = 20
$ws.Cells.Item(19,8).Value = 'https://show.bilibili.com/platform/detail.html?id=89200'
$ws.Cells.Item(19,9).Value = '//i2.hdslb.com/bfs/openplatform/202407/b4w7Ifkm1720766324652.jpeg'

$ws.Cells.Item(20,2).Value = '2024-07-27'
$ws.Cells.Item(20,3).Value = '上海·创世次元动漫游戏嘉年华6.0'
$ws.Cells.Item(20,4).Value = '老沪闵路1388号舒也时代广场C栋2层 轮客行轮滑馆(闵行店)'
$ws.Cells.Item(20,5).Value = '2024.07.27 10:00-07.28 17:00'
$ws.Cells.Item(20,6).Value = 1212
$ws.Cells.Item(20,7).Value = 65
$ws.Cells.Item(20,8).Value = 'https://show.bilibili.com/platform/detail.html?id=89043'
$ws.Cells.Item(20,9).Value = '//i1.hdslb.com/bfs/openplatform/202407/rBCZUYJ61720596521825.jpeg'

$ws.Cells.Item(21,2).Value = '2024-07-27'
$ws.Cells.Item(21,3).Value = '上海·名侦探柯南 连载30周年纪念展（早鸟票）'
$ws.Cells.Item(21,4).Value = '南京西路2-68号 新世界城11楼'
$ws.Cells.Item(21,5).Value = '2024.07.27 00:00-08.04 23:59'
$ws.Cells.Item(21,6).Value = 464
$ws.Cells.Item(21,7).Value = 79
$ws.Cells.Item(21,8).Value = 'https://show.bilibili.com/platform/detail.html?id=89294'
$ws.Cells.Item(21,9).Value = '//i0.hdslb.com/bfs/openplatform/202407/35thNBrO1721035918311.png'
